$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (B1:J1) - columns permuted (cyclic shift of the 9 field columns)
$ws.Range("B1").Value2 = "way_of_speech"
$ws.Range("C1").Value2 = "socio_economic"
$ws.Range("D1").Value2 = "ethnicity_skin_color"
$ws.Range("E1").Value2 = "personality"
$ws.Range("F1").Value2 = "dress_propeties"
$ws.Range("G1").Value2 = "political_affiliation"
$ws.Range("H1").Value2 = "hobbies"
$ws.Range("I1").Value2 = "body_size"
$ws.Range("J1").Value2 = "intelligence"

# Row 2 (B2:J2) data permuted to match the new column order
$ws.Range("B2").Value2 = 0.153503587410487
$ws.Range("C2").Value2 = 0.0659424687055382
$ws.Range("D2").Value2 = 0.0608644296500644
$ws.Range("E2").Value2 = 0.371233119411251
$ws.Range("F2").Value2 = 0.114690402043858
$ws.Range("G2").Value2 = 0.0246496399176212
$ws.Range("H2").Value2 = 0.0246354565925125
$ws.Range("I2").Value2 = 0.0639616461422115
$ws.Range("J2").Value2 = 0.120519250126453

# Row 3 (B3:J3) data permuted to match the new column order
$ws.Range("B3").Value2 = 0.049367251588039
$ws.Range("C3").Value2 = 0.0781807725535249
$ws.Range("D3").Value2 = 0.08035278368559
$ws.Range("E3").Value2 = 0.264966338380914
$ws.Range("F3").Value2 = 0.0761646517897155
$ws.Range("G3").Value2 = 0.0500482628634886
$ws.Range("H3").Value2 = 0.124465299991199
$ws.Range("I3").Value2 = 0.141890646193694
$ws.Range("J3").Value2 = 0.134563992953834

# Row 4 (B4:J4) data permuted to match the new column order
$ws.Range("B4").Value2 = 0.0612225992848502
$ws.Range("C4").Value2 = 0.0959631394816918
$ws.Range("D4").Value2 = 0.0816216309435467
$ws.Range("E4").Value2 = 0.217040717526421
$ws.Range("F4").Value2 = 0.055746416527572
$ws.Range("G4").Value2 = 0.0782812153499518
$ws.Range("H4").Value2 = 0.245964213473663
$ws.Range("I4").Value2 = 0.0666572745306769
$ws.Range("J4").Value2 = 0.0975027928816254

# Row 5 (B5:J5) data permuted to match the new column order
$ws.Range("B5").Value2 = 0.110341930877243
$ws.Range("C5").Value2 = 0.0344324152668377
$ws.Range("D5").Value2 = 0.11173760988235
$ws.Range("E5").Value2 = 0.151436700089937
$ws.Range("F5").Value2 = 0.116615554321736
$ws.Range("G5").Value2 = 0.0393104701258311
$ws.Range("H5").Value2 = 0.0570544782861202
$ws.Range("I5").Value2 = 0.349001313211675
$ws.Range("J5").Value2 = 0.0300695279382673

# Row 6 (B6:J6, T6:AB6) - replaced with newly collected participant data
$ws.Range("B6").Value2 = -0.014
$ws.Range("C6").Value2 = 0.119
$ws.Range("D6").Value2 = 0.089
$ws.Range("E6").Value2 = 0.144
$ws.Range("F6").Value2 = 0.073
$ws.Range("G6").Value2 = 0.254
$ws.Range("H6").Value2 = 0.115
$ws.Range("I6").Value2 = 0.118
$ws.Range("J6").Value2 = 0.136
$ws.Range("T6").Value2 = 0.086
$ws.Range("U6").Value2 = 0.116
$ws.Range("V6").Value2 = 0.195
$ws.Range("W6").Value2 = 0.124
$ws.Range("X6").Value2 = 0.055
$ws.Range("Y6").Value2 = 0.139
$ws.Range("Z6").Value2 = 0.046
$ws.Range("AA6").Value2 = 0.128
$ws.Range("AB6").Value2 = 0.178

# Update the active selection to match the saved view state
$ws.Range("I6").Select()
